# Negate every numeric value in column "Block" (column E) across the
# data rows of the frame-data sheet. This flips the sign convention for
# block frame advantage (e.g. 2 -> -2, -6 -> 6), leaving blank/empty
# cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [double]) {
        if ($val -ne 0) {
            $cell.Value2 = -$val
        }
    }
}
